$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 500.1
$ws.Range("J17").Value = 500.1
$ws.Range("L17").Value = 1500.3
$ws.Range("N17").Value = -1836.3

$ws.Range("H137").Value = 1538.1666
$ws.Range("I137").Value = 1015.1429
$ws.Range("J137").Value = 1871
$ws.Range("K137").Value = 3045.4287
$ws.Range("L137").Value = 5613
$ws.Range("M137").Value = -495.4287000000004
$ws.Range("N137").Value = -10713

$ws.Range("H138").Value = 1397.798
$ws.Range("J138").Value = 1995.9298
$ws.Range("L138").Value = 5987.7894
$ws.Range("N138").Value = -16267.7894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1080.0588
$ws.Range("I45").Value = 1081.6923
$ws.Range("J45").Value = 1074.75
$ws.Range("K45").Value = 1081.6923
$ws.Range("L45").Value = 1074.75
$ws.Range("M45").Value = -704.6922999999999
$ws.Range("N45").Value = -1828.75

$ws.Range("H61").Value = 1020.9667
$ws.Range("I61").Value = 993.6667
$ws.Range("J61").Value = 1266.6666
$ws.Range("K61").Value = 993.6667
$ws.Range("L61").Value = 1266.6666
$ws.Range("M61").Value = -781.6667
$ws.Range("N61").Value = -1690.6666

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H136").Value = 1020.9667
$ws.Range("I136").Value = 993.6667
$ws.Range("J136").Value = 1266.6666
$ws.Range("K136").Value = 2981.0001
$ws.Range("L136").Value = 3799.9998
$ws.Range("M136").Value = -431.0001000000002
$ws.Range("N136").Value = -8899.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2844.41
$ws.Range("I134").Value = 816.69385
$ws.Range("J134").Value = 11124.25
$ws.Range("K134").Value = 2450.08155
$ws.Range("L134").Value = 33372.75
$ws.Range("M134").Value = 84.91845000000012
$ws.Range("N134").Value = -38442.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1195.6666
$ws.Range("I31").Value = 1195.6666
$ws.Range("K31").Value = 1195.6666
$ws.Range("M31").Value = -900.6666

$ws.Range("H34").Value = 1195.6666
$ws.Range("I34").Value = 1195.6666
$ws.Range("K34").Value = 1195.6666
$ws.Range("M34").Value = -993.6666

$ws.Range("H58").Value = 1772.5
$ws.Range("I58").Value = 1580
$ws.Range("J58").Value = 2350
$ws.Range("K58").Value = 1580
$ws.Range("L58").Value = 2350
$ws.Range("M58").Value = -1377
$ws.Range("N58").Value = -2756

$ws.Range("H132").Value = 3799.6667
$ws.Range("I132").Value = 1800
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 5400
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -2870
$ws.Range("N132").Value = -17658.8

$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -40060

$ws.Range("H136").Value = 1772.5
$ws.Range("I136").Value = 1580
$ws.Range("J136").Value = 2350
$ws.Range("K136").Value = 4740
$ws.Range("L136").Value = 7050
$ws.Range("M136").Value = -2190
$ws.Range("N136").Value = -12150

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 74
$ws.Range("I12").Value = 99.333336
$ws.Range("J12").Value = 66
$ws.Range("K12").Value = 298.000008
$ws.Range("L12").Value = 198
$ws.Range("M12").Value = -125.000008
$ws.Range("N12").Value = -544

$ws.Range("H131").Value = 14495145
$ws.Range("J131").Value = 2645.541
$ws.Range("L131").Value = 7936.623000000001
$ws.Range("N131").Value = -18016.623

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2082
$ws.Range("I122").Value = 2526.7144
$ws.Range("J122").Value = 1563.1666
$ws.Range("K122").Value = 7580.1432
$ws.Range("L122").Value = 4689.4998
$ws.Range("M122").Value = -5130.1432
$ws.Range("N122").Value = -9589.4998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2940.1333
$ws.Range("I40").Value = 2866.1667
$ws.Range("J40").Value = 2989.4443
$ws.Range("K40").Value = 2866.1667
$ws.Range("L40").Value = 2989.4443
$ws.Range("M40").Value = -2730.1667
$ws.Range("N40").Value = -3261.4443

$ws.Range("H61").Value = 1733.1111
$ws.Range("I61").Value = 1448
$ws.Range("K61").Value = 1448
$ws.Range("M61").Value = -1246

$ws.Range("H68").Value = 1921
$ws.Range("I68").Value = 1802
$ws.Range("J68").Value = 2099.5
$ws.Range("K68").Value = 1802
$ws.Range("L68").Value = 2099.5
$ws.Range("M68").Value = -1053
$ws.Range("N68").Value = -3597.5

$ws.Range("H71").Value = 1921
$ws.Range("I71").Value = 1802
$ws.Range("J71").Value = 2099.5
$ws.Range("K71").Value = 9010
$ws.Range("L71").Value = 10497.5
$ws.Range("M71").Value = -5266
$ws.Range("N71").Value = -17985.5

$ws.Range("H82").Value = 1405.6207
$ws.Range("I82").Value = 1320.1052
$ws.Range("J82").Value = 1568.1
$ws.Range("K82").Value = 1320.1052
$ws.Range("L82").Value = 1568.1
$ws.Range("M82").Value = -959.1052
$ws.Range("N82").Value = -2290.1

$ws.Range("H85").Value = 1405.6207
$ws.Range("I85").Value = 1320.1052
$ws.Range("J85").Value = 1568.1
$ws.Range("K85").Value = 1320.1052
$ws.Range("L85").Value = 1568.1
$ws.Range("M85").Value = -72.10519999999997
$ws.Range("N85").Value = -4064.1

$ws.Range("H113").Value = 1733.1111
$ws.Range("I113").Value = 1448
$ws.Range("K113").Value = 1448
$ws.Range("M113").Value = 722

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 83338170
$ws.Range("I62").Value = 100005200
$ws.Range("J62").Value = 3003
$ws.Range("K62").Value = 100005200
$ws.Range("L62").Value = 3003
$ws.Range("M62").Value = -100004576
$ws.Range("N62").Value = -4251

$ws.Range("H65").Value = 83338170
$ws.Range("I65").Value = 100005200
$ws.Range("J65").Value = 3003
$ws.Range("K65").Value = 500026000
$ws.Range("L65").Value = 15015
$ws.Range("M65").Value = -500022880
$ws.Range("N65").Value = -21255

$ws.Range("H107").Value = 640
$ws.Range("J107").Value = 745
$ws.Range("L107").Value = 2235
$ws.Range("N107").Value = -6075

$ws.Range("H122").Value = 96924670
$ws.Range("I122").Value = 126001540
$ws.Range("K122").Value = 378004620
$ws.Range("M122").Value = -378002170
